$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet stored a project-specific data-lab access code ("[DL-MAA2020-01]")
# in column C (rows 10-13). Replace it with the generic placeholder code
# "[DL-MAA20XX-YY]" so the workbook can be shared as a template.
$ws.Range("C10:C13").Value = "[DL-MAA20XX-YY]"
